# Adds the Dragon Boat Festival entries to the "diary":
#   - appends "2022年6月3日星期五" into the (previously empty) paragraph
#     that follows the "...明天就是端午节了" entry (the one holding the
#     _GoBack bookmark);
#   - inserts a brand new paragraph right after it with the festival-day
#     weather note, carrying the _GoBack bookmark forward onto it;
#   - inserts a new, empty paragraph after that one.
#
# Both new runs need the same "typed with an East-Asian IME" formatting
# ( <w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/> )
# as every other CJK run already in the document, and the new paragraph
# marks need <w:rFonts w:hint="default"/><w:lang .../> to match their
# siblings. Range.InsertXML lets us stamp that markup exactly instead of
# hoping the object model infers it from plain text.

$d = $word.ActiveDocument

# Locate the paragraph that used to be empty (it only held the
# bookmarkStart/bookmarkEnd for "_GoBack") -- it is the paragraph right
# after the one that announces "明天就是端午节了" (tomorrow is the
# Dragon Boat Festival). (Paragraph.Next is unreliable in this host, so
# find the index via the collection and re-fetch by Item() instead.)
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*明天就是端午节了*") {
        $anchorIndex = $i + 1
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate the anchor paragraph"
}
$anchor = $d.Paragraphs.Item($anchorIndex)
if ($anchor -eq $null) {
    throw "Could not locate the target (bookmark) paragraph"
}

$insertPoint = $d.Range($anchor.Range.Start, $anchor.Range.Start)

# First payload: a new paragraph holding the date line, followed by a
# second paragraph holding the festival weather note -- the trailing
# paragraph of the inserted XML merges into the destination paragraph,
# so the bookmark that already lives there rides along onto the
# "festival note" paragraph, exactly like the diff shows.
$datePara = '<w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
    '<w:t>2022年6月3日星期五</w:t></w:r>' +
    '</w:p>'
$festivalPara = '<w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
    '<w:t>中雨，今天是农历五月初五，中国传统端午节</w:t></w:r>' +
    '</w:p>'

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $datePara + $festivalPara + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$insertPoint.InsertXML($xml1)

# Second payload: a brand new empty paragraph inserted right after the
# festival-note paragraph (which now carries the bookmark).
$festival = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*中国传统端午节*") {
        $festival = $p
        break
    }
}
if ($festival -eq $null) {
    throw "Could not locate the newly inserted festival paragraph"
}

$afterPoint = $d.Range($festival.Range.End, $festival.Range.End)
$emptyPara = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>'
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $emptyPara + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$afterPoint.InsertXML($xml2)

Write-Output "Dragon Boat Festival entries inserted."
